$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "66.009.88"
$ws.Range("E2").Value = "  -2.07%  "
Set-TextValue $ws.Range("D3") "3.446.34"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "583.45"
$ws.Range("E5").Value = "  -1.70%  "
Set-TextValue $ws.Range("D6") "173.30"
$ws.Range("E6").Value = "  -4.33%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.595"
$ws.Range("E8").Value = "  -3.63%  "
Set-TextValue $ws.Range("D9") "3.441.67"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("E10").Value = "  -6.48%  "
Set-TextValue $ws.Range("D11") "6.86"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("E12").Value = "  -4.48%  "
Set-TextValue $ws.Range("D13") "4.036.94"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  -0.29%  "
Set-TextValue $ws.Range("D15") "29.94"
$ws.Range("E15").Value = "  -6.41%  "
Set-TextValue $ws.Range("D16") "66.062.15"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("E17").Value = "  -3.72%  "
Set-TextValue $ws.Range("D18") "3.440.52"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("E19").Value = "  -4.98%  "
Set-TextValue $ws.Range("D20") "13.79"
$ws.Range("E20").Value = "  -2.22%  "
Set-TextValue $ws.Range("D21") "366.24"
$ws.Range("E21").Value = "  -7.06%  "
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue $ws.Range("D24") "72.02"
$ws.Range("E24").Value = "  +0.31%  "
Set-TextValue $ws.Range("D25") "0.529"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("E26").Value = "  -3.02%  "
Set-TextValue $ws.Range("D27") "9.65"
$ws.Range("E27").Value = "  -7.04%  "
$ws.Range("E28").Value = "  +1.09%  "
Set-TextValue $ws.Range("D29") "0.997"
$ws.Range("E29").Value = "  -0.32%  "
Set-TextValue $ws.Range("D30") "23.99"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.98"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "5.74"
$ws.Range("E32").Value = "  -6.47%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -8.31%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  -2.59%  "
Set-TextValue $ws.Range("D37") "157.79"
$ws.Range("E37").Value = "  -2.00%  "
Set-TextValue $ws.Range("D38") "29.05"
$ws.Range("E38").Value = "  +10.69%  "
$ws.Range("E39").Value = "  -1.55%  "
$ws.Range("E40").Value = "  -4.51%  "
Set-TextValue $ws.Range("D41") "2.56"
$ws.Range("E41").Value = "  -10.10%  "
Set-TextValue $ws.Range("D42") "2.739.31"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -5.86%  "
Set-TextValue $ws.Range("D44") "6.28"
$ws.Range("E44").Value = "  -7.33%  "
Set-TextValue $ws.Range("D45") "0.0683"
$ws.Range("E45").Value = "  -4.82%  "
Set-TextValue $ws.Range("D46") "39.91"
$ws.Range("E46").Value = "  -3.89%  "
Set-TextValue $ws.Range("D47") "24.16"
$ws.Range("E47").Value = "  -8.10%  "
Set-TextValue $ws.Range("D48") "0.0287"
$ws.Range("E48").Value = "  -3.77%  "
Set-TextValue $ws.Range("D49") "303.98"
$ws.Range("E49").Value = "  -6.77%  "
Set-TextValue $ws.Range("D50") "0.814"
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("E51").Value = "  -3.67%  "

Write-Output "Applied cryptos list update"
